$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for rows 6-14: Column C (Minimum Order Quantity), Column E (Final Price), Column F (Cost Price)
# Column B (NSN Number / "Rnd Manufacturer") is cleared for all these rows.
$data = @(
    @{ Row = 6;  C = 1; E = 112.17; F = 7.3 },
    @{ Row = 7;  C = 1; E = 200;    F = 4.34 },
    @{ Row = 8;  C = 1; E = 114.2;  F = 74.23 },
    @{ Row = 9;  C = 2; E = 27.67;  F = 18.56 },
    @{ Row = 10; C = 1; E = 505;    F = 411.14 },
    @{ Row = 11; C = 1; E = 800;    F = 118.84 },
    @{ Row = 12; C = 1; E = 32.25;  F = 20.96 },
    @{ Row = 13; C = 1; E = 812.01; F = 682.09 },
    @{ Row = 14; C = 1; E = 28.68;  F = 1.29 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).ClearContents()
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
}
